$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 557.6
$ws.Range("I4").Value = 94.5
$ws.Range("J4").Value = 866.3333
$ws.Range("K4").Value = 94.5
$ws.Range("L4").Value = 866.3333
$ws.Range("M4").Value = 19.5
$ws.Range("N4").Value = -1094.3333
$ws.Range("H98").Value = 1604.4
$ws.Range("J98").Value = 3566.6667
$ws.Range("L98").Value = 3566.6667
$ws.Range("N98").Value = -6562.6667
$ws.Range("H101").Value = 217.57143
$ws.Range("I101").Value = 210.75
$ws.Range("J101").Value = 226.66667
$ws.Range("K101").Value = 632.25
$ws.Range("L101").Value = 680.00001
$ws.Range("M101").Value = 989.75
$ws.Range("N101").Value = -3924.00001
$ws.Range("H122").Value = 1604.4
$ws.Range("J122").Value = 3566.6667
$ws.Range("L122").Value = 10700.0001
$ws.Range("N122").Value = -15600.0001
$ws.Range("H125").Value = 3020.75
$ws.Range("I125").Value = 2466.6667
$ws.Range("K125").Value = 22200.0003
$ws.Range("M125").Value = -19740.0003
$ws.Range("H132").Value = 16484.559
$ws.Range("I132").Value = 1817.5084
$ws.Range("K132").Value = 5452.5252
$ws.Range("M132").Value = -2922.5252
$ws.Range("H135").Value = 913.8333
$ws.Range("I135").Value = 913.8333
$ws.Range("K135").Value = 8224.4997
$ws.Range("M135").Value = -5689.4997
$ws.Range("H138").Value = 8775473
$ws.Range("I138").Value = 2085.125
$ws.Range("K138").Value = 6255.375
$ws.Range("M138").Value = -1115.375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 573.625
$ws.Range("I2").Value = 404.72726
$ws.Range("K2").Value = 404.72726
$ws.Range("M2").Value = -291.72726
$ws.Range("H32").Value = 6644.036
$ws.Range("I32").Value = 3793.6025
$ws.Range("K32").Value = 3793.6025
$ws.Range("M32").Value = -3506.6025
$ws.Range("H33").Value = 21342
$ws.Range("I33").Value = 21342
$ws.Range("K33").Value = 21342
$ws.Range("M33").Value = -21013
$ws.Range("H61").Value = 405086.7
$ws.Range("I61").Value = 3060
$ws.Range("K61").Value = 3060
$ws.Range("M61").Value = -2848
$ws.Range("H116").Value = 573.625
$ws.Range("I116").Value = 404.72726
$ws.Range("K116").Value = 404.72726
$ws.Range("M116").Value = 1889.27274
$ws.Range("H122").Value = 1226.5714
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H132").Value = 3551.682
$ws.Range("I132").Value = 3396.9268
$ws.Range("K132").Value = 10190.7804
$ws.Range("M132").Value = -7660.7804
$ws.Range("H136").Value = 405086.7
$ws.Range("I136").Value = 3060
$ws.Range("K136").Value = 9180
$ws.Range("M136").Value = -6630
$ws.Range("N122").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 573.625
$ws.Range("I3").Value = 404.72726
$ws.Range("K3").Value = 404.72726
$ws.Range("M3").Value = -290.72726
$ws.Range("H20").Value = 5059
$ws.Range("I20").Value = 4694.8887
$ws.Range("J20").Value = 5787.222
$ws.Range("K20").Value = 4694.8887
$ws.Range("L20").Value = 5787.222
$ws.Range("M20").Value = -4447.8887
$ws.Range("N20").Value = -6281.222
$ws.Range("H80").Value = 676.6875
$ws.Range("I80").Value = 510.7143
$ws.Range("K80").Value = 510.7143
$ws.Range("M80").Value = 487.2857
$ws.Range("H83").Value = 676.6875
$ws.Range("I83").Value = 510.7143
$ws.Range("K83").Value = 2553.5715
$ws.Range("M83").Value = 2438.4285
$ws.Range("H134").Value = 1896.3062
$ws.Range("I134").Value = 1856.0652
$ws.Range("K134").Value = 5568.1956
$ws.Range("M134").Value = -3033.1956

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 49437.855
$ws.Range("I31").Value = 57028.445
$ws.Range("K31").Value = 57028.445
$ws.Range("M31").Value = -56733.445
$ws.Range("H34").Value = 49437.855
$ws.Range("I34").Value = 57028.445
$ws.Range("K34").Value = 57028.445
$ws.Range("M34").Value = -56826.445
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H127").Value = 141793.75
$ws.Range("J127").Value = 141793.75
$ws.Range("L127").Value = 141793.75
$ws.Range("N127").Value = -151713.75
$ws.Range("H130").Value = 80999.5
$ws.Range("J130").Value = 80999.5
$ws.Range("L130").Value = 80999.5
$ws.Range("N130").Value = -91039.5
$ws.Range("H132").Value = 1762.091
$ws.Range("I132").Value = 1592.3158
$ws.Range("J132").Value = 2837.3333
$ws.Range("K132").Value = 4776.9474
$ws.Range("L132").Value = 8511.999899999999
$ws.Range("M132").Value = -2246.9474
$ws.Range("N132").Value = -13571.9999
$ws.Range("H133").Value = 73500
$ws.Range("J133").Value = 73500
$ws.Range("L133").Value = 73500
$ws.Range("N133").Value = -78560
$ws.Range("H141").Value = 204107.6
$ws.Range("J141").Value = 204107.6
$ws.Range("L141").Value = 204107.6
$ws.Range("N141").Value = -214467.6
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 55555660
$ws.Range("I18").Value = 55555660
$ws.Range("K18").Value = 55555660
$ws.Range("M18").Value = -55555367
$ws.Range("H55").Value = 21375
$ws.Range("I55").Value = 12000
$ws.Range("J55").Value = 24500
$ws.Range("K55").Value = 12000
$ws.Range("L55").Value = 24500
$ws.Range("M55").Value = -11673
$ws.Range("N55").Value = -25154
$ws.Range("H70").Value = 10794.115
$ws.Range("I70").Value = 11685.23
$ws.Range("J70").Value = 9903
$ws.Range("K70").Value = 11685.23
$ws.Range("L70").Value = 9903
$ws.Range("M70").Value = -11415.23
$ws.Range("N70").Value = -10443
$ws.Range("H73").Value = 10794.115
$ws.Range("I73").Value = 11685.23
$ws.Range("J73").Value = 9903
$ws.Range("K73").Value = 11685.23
$ws.Range("L73").Value = 9903
$ws.Range("M73").Value = -10749.23
$ws.Range("N73").Value = -11775
$ws.Range("H122").Value = 2427.6667
$ws.Range("J122").Value = 2999
$ws.Range("L122").Value = 8997
$ws.Range("N122").Value = -13897
$ws.Range("H136").Value = 43099.4
$ws.Range("J136").Value = 43099.4
$ws.Range("L136").Value = 129298.2
$ws.Range("N136").Value = -134398.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3505.5
$ws.Range("I61").Value = 2941.2666
$ws.Range("K61").Value = 2941.2666
$ws.Range("M61").Value = -2739.2666
$ws.Range("H68").Value = 2250.95
$ws.Range("I68").Value = 2378.2144
$ws.Range("J68").Value = 1954
$ws.Range("K68").Value = 2378.2144
$ws.Range("L68").Value = 1954
$ws.Range("M68").Value = -1629.2144
$ws.Range("N68").Value = -3452
$ws.Range("H71").Value = 2250.95
$ws.Range("I71").Value = 2378.2144
$ws.Range("J71").Value = 1954
$ws.Range("K71").Value = 11891.072
$ws.Range("L71").Value = 9770
$ws.Range("M71").Value = -8147.072
$ws.Range("N71").Value = -17258
$ws.Range("H113").Value = 3505.5
$ws.Range("I113").Value = 2941.2666
$ws.Range("K113").Value = 2941.2666
$ws.Range("M113").Value = -771.2665999999999
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 26949
$ws.Range("J130").Value = 26949
$ws.Range("L130").Value = 26949
$ws.Range("N130").Value = -36989
$ws.Range("H131").Value = 154992.5
$ws.Range("J131").Value = 154992.5
$ws.Range("L131").Value = 154992.5
